# Daily attendance processing - 2025-10-23 22:20:57
# Normalizes the "Recorded By" (column G) text on the "Session Analysis Results"
# sheet: reorders the comma-separated recorder names.
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "backup@backdoor.com, system, System"     -> "backup@backdoor.com, System, system"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$rows = @(2,3,6,10,11,12,13,14,15,17,18,19,20,21,22,29,30,33,37,38,39,40,41,42,44,45,46,47,48,49,56,57,60,64,65,66,67,68,69,71,72,73,74,75,76,86,87,88,89,93,95,96,97,99,102,112,113,114,115,119,121,122,123,125,128,138,139,140,141,145,147,148,149,151,154)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Value2
    if ($value -eq "dnasr281@gmail.com, System") {
        $cell.Value = "System, dnasr281@gmail.com"
    } elseif ($value -eq "backup@backdoor.com, system, System") {
        $cell.Value = "backup@backdoor.com, System, system"
    } else {
        # Fallback: generically swap the last two comma-separated parts
        $parts = $value -split ", "
        if ($parts.Length -ge 2) {
            $last = $parts[$parts.Length - 1]
            $secondLast = $parts[$parts.Length - 2]
            $parts[$parts.Length - 1] = $secondLast
            $parts[$parts.Length - 2] = $last
            $cell.Value = [string]::Join(", ", $parts)
        }
    }
}
